$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the canvas_role test data values (F3/F4) to the roles the test now
# expects, per the commit message ("Updated the test data to reflect the
# canvas role the test is looking for").
$ws.Range("F3").Value = "Observer"
$ws.Range("F4").Value = "TA"

# Resize a few columns (canvas_role related columns) so their contents are
# fully visible, and leave the cursor/selection on A2 as it was left in the
# saved workbook.
$ws.Columns.Item(3).ColumnWidth = 21.1
$ws.Columns.Item(4).ColumnWidth = 16.8
$ws.Columns.Item(6).ColumnWidth = 25.3

[void]$ws.Range("A2").Select()
